# Apply scheduled-runner updates to the Mateus_Profits leve profitability sheets.
# Generated from the authoritative cell-level diff against before.xlsx.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -797
$ws.Range("N34").ClearContents() | Out-Null

$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -285
$ws.Range("N36").ClearContents() | Out-Null

$ws.Range("H40").Value = 5051.4
$ws.Range("I40").Value = 2922.8
$ws.Range("K40").Value = 2922.8
$ws.Range("M40").Value = -2747.8

$ws.Range("H43").Value = 5981.8335
$ws.Range("J43").Value = 6178.2
$ws.Range("L43").Value = 6178.2
$ws.Range("N43").Value = -6316.2

$ws.Range("H51").Value = 9949.85
$ws.Range("J51").Value = 10088.059
$ws.Range("L51").Value = 10088.059
$ws.Range("N51").Value = -11056.059

$ws.Range("H92").Value = 412.6
$ws.Range("I92").Value = 412.6
$ws.Range("K92").Value = 412.6
$ws.Range("M92").Value = 835.4

$ws.Range("H105").Value = 63650
$ws.Range("J105").Value = 63650
$ws.Range("L105").Value = 63650
$ws.Range("N105").Value = -70638

$ws.Range("H112").Value = 3638.6
$ws.Range("J112").Value = 3638.6
$ws.Range("L112").Value = 10915.8
$ws.Range("N112").Value = -13131.8

$ws.Range("H113").Value = 125004360
$ws.Range("I113").Value = 166669300
$ws.Range("K113").Value = 166669300
$ws.Range("M113").Value = -166666046

$ws.Range("H131").Value = 5662.8335
$ws.Range("I131").Value = 5595.4
$ws.Range("J131").Value = 6000
$ws.Range("K131").Value = 16786.2
$ws.Range("L131").Value = 18000
$ws.Range("M131").Value = -11746.2
$ws.Range("N131").Value = -28080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1594

$ws.Range("H37").Value = 19200
$ws.Range("J37").Value = 21000
$ws.Range("L37").Value = 21000
$ws.Range("N37").Value = -21546

$ws.Range("H46").Value = 7484
$ws.Range("J46").Value = 8420
$ws.Range("L46").Value = 8420
$ws.Range("N46").Value = -9058

$ws.Range("H88").Value = 3112.8333
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 3235.4
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 3235.4
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -4047.4

$ws.Range("H91").Value = 3112.8333
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 3235.4
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 3235.4
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -6043.4

$ws.Range("H122").Value = 2147.1667
$ws.Range("I122").Value = 1685.3334
$ws.Range("J122").Value = 3532.6667
$ws.Range("K122").Value = 5056.0002
$ws.Range("L122").Value = 10598.0001
$ws.Range("M122").Value = -2606.0002
$ws.Range("N122").Value = -15498.0001

$ws.Range("H132").Value = 4257.6045
$ws.Range("I132").Value = 3513.6
$ws.Range("K132").Value = 10540.8
$ws.Range("M132").Value = -8010.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 12396.6
$ws.Range("I36").Value = 10495.75
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 10495.75
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -9961.75
$ws.Range("N36").Value = -21068

$ws.Range("H86").Value = 2696.111
$ws.Range("I86").Value = 2568.077
$ws.Range("J86").Value = 3029
$ws.Range("K86").Value = 2568.077
$ws.Range("L86").Value = 3029
$ws.Range("M86").Value = -1445.077
$ws.Range("N86").Value = -5275

$ws.Range("H89").Value = 2696.111
$ws.Range("I89").Value = 2568.077
$ws.Range("J89").Value = 3029
$ws.Range("K89").Value = 12840.385
$ws.Range("L89").Value = 15145
$ws.Range("M89").Value = -7224.385000000002
$ws.Range("N89").Value = -26377

$ws.Range("H106").Value = 14999.5
$ws.Range("J106").Value = 14999.5
$ws.Range("L106").Value = 14999.5
$ws.Range("N106").Value = -17523.5

$ws.Range("H134").Value = 4247.871
$ws.Range("I134").Value = 4247.871
$ws.Range("K134").Value = 12743.613
$ws.Range("M134").Value = -10208.613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4746.8335
$ws.Range("I16").Value = 3665.3333
$ws.Range("K16").Value = 3665.3333
$ws.Range("M16").Value = -3378.3333

$ws.Range("H95").Value = 42333.332
$ws.Range("J95").Value = 42333.332
$ws.Range("L95").Value = 42333.332
$ws.Range("N95").Value = -47825.332

$ws.Range("H105").Value = 2368.6365
$ws.Range("I105").Value = 2368.6365
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2368.6365
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -621.6365000000001
$ws.Range("N105").ClearContents() | Out-Null

$ws.Range("H113").Value = 4746.8335
$ws.Range("I113").Value = 3665.3333
$ws.Range("K113").Value = 3665.3333
$ws.Range("M113").Value = -1495.3333

$ws.Range("H122").Value = 3188.5652
$ws.Range("I122").Value = 2633.5264
$ws.Range("J122").Value = 5825
$ws.Range("K122").Value = 7900.5792
$ws.Range("L122").Value = 17475
$ws.Range("M122").Value = -5450.5792
$ws.Range("N122").Value = -22375

$ws.Range("H132").Value = 2271.842
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 1300
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 3900
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -3089
$ws.Range("N68").Value = -6122

$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 1300
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 11700
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -7644
$ws.Range("N71").Value = -21612

$ws.Range("H97").Value = 3024.7273
$ws.Range("I97").Value = 3480.8333
$ws.Range("J97").Value = 2477.4
$ws.Range("K97").Value = 10442.4999
$ws.Range("L97").Value = 7432.200000000001
$ws.Range("M97").Value = -9946.499899999999
$ws.Range("N97").Value = -8424.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 15000
$ws.Range("K31").Value = 15000
$ws.Range("M31").Value = -14708

$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 15000
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = -14723

$ws.Range("H97").Value = 3853.353
$ws.Range("I97").Value = 731.46155
$ws.Range("K97").Value = 731.46155
$ws.Range("M97").Value = -235.46155

$ws.Range("H122").Value = 3075.7368
$ws.Range("I122").Value = 2468.125
$ws.Range("J122").Value = 6316.3335
$ws.Range("K122").Value = 7404.375
$ws.Range("L122").Value = 18949.0005
$ws.Range("M122").Value = -4954.375
$ws.Range("N122").Value = -23849.0005

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 15483.75
$ws.Range("I42").Value = 15999
$ws.Range("J42").Value = 14968.5
$ws.Range("K42").Value = 15999
$ws.Range("L42").Value = 14968.5
$ws.Range("M42").Value = -15436
$ws.Range("N42").Value = -16094.5

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents() | Out-Null

$ws.Range("H49").Value = 15483.75
$ws.Range("I49").Value = 15999
$ws.Range("J49").Value = 14968.5
$ws.Range("K49").Value = 15999
$ws.Range("L49").Value = 14968.5
$ws.Range("M49").Value = -15852
$ws.Range("N49").Value = -15262.5

$ws.Range("H93").Value = 2450
$ws.Range("I93").Value = 2425
$ws.Range("K93").Value = 2425
$ws.Range("M93").Value = -1177

$ws.Range("H132").Value = 7274.9575
$ws.Range("I132").Value = 7003.1953
$ws.Range("J132").Value = 9132
$ws.Range("K132").Value = 21009.5859
$ws.Range("L132").Value = 27396
$ws.Range("M132").Value = -18479.5859
$ws.Range("N132").Value = -32456

$ws.Range("H135").Value = 61357
$ws.Range("J135").Value = 61357
$ws.Range("L135").Value = 61357
$ws.Range("N135").Value = -71497

$ws.Range("H140").Value = 95590
$ws.Range("J140").Value = 95590
$ws.Range("L140").Value = 95590
$ws.Range("N140").Value = -105950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").ClearContents() | Out-Null

$ws.Range("H23").Value = 4997.3335
$ws.Range("I23").Value = 4992
$ws.Range("K23").Value = 4992
$ws.Range("M23").Value = -4763

$ws.Range("H70").Value = 38277.145
$ws.Range("J70").Value = 40490
$ws.Range("L70").Value = 40490
$ws.Range("N70").Value = -41120

$ws.Range("H73").Value = 38277.145
$ws.Range("J73").Value = 40490
$ws.Range("L73").Value = 40490
$ws.Range("N73").Value = -42674

$ws.Range("H96").Value = 3922
$ws.Range("I96").Value = 4383
$ws.Range("K96").Value = 4383
$ws.Range("M96").Value = -3010

$ws.Range("H103").Value = 48490
$ws.Range("J103").Value = 48490
$ws.Range("L103").Value = 48490
$ws.Range("N103").Value = -50834
